$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.691.61'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.80%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.635.24'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.88%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.46'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.10%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.493'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.71%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.41%  '

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.81%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.02'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +4.44%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0837'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.91%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.865.74'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +2.01%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.635.07'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.86%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.06'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.25%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.525'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.48%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.709.65'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.95%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.00'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.89%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0₃0740'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.84%  '

$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.00'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.06%  '

$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '208.32'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +3.75%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.31'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.85%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.39'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.17%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.13'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.27%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.91'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +2.68%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.65%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.06%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.53%  '

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +3.16%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.39'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.28%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0515'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +5.45%  '

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.33%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.77%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.96'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.28%  '

$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.42'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.11%  '

$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = 'LidoDAOToken'
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.50'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.82%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.167.30'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.13%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.19%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.809'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +3.08%  '

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.01%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.503'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.71%  '

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.22%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.793'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.76%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.37'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.53%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.775.26'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +2.01%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.47'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.00%  '

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.40%  '

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +10.12%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '54.91'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.57%  '

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.38%  '

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.69%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.52'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +4.08%  '
